$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5301_topic_13__ID** " -> "**ID__AFFARS_5301_402__ID**" ---
$p = $d.Paragraphs(1)

$oldId = "**ID__AFFARS_5301_topic_13__ID**"
$newId = "**ID__AFFARS_5301_402__ID**"

# Locate the placeholder text (it lives in the paragraph's first run) without
# disturbing the separate run that holds the trailing space.
$findRange = $p.Range.Duplicate
$findRange.Find.Execute($oldId, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($findRange.Find.Found) {
    $matchStart = $findRange.Start
    $matchEnd = $findRange.End

    $target = $d.Range($matchStart, $matchEnd)
    $target.Text = $newId

    # Drop whatever text remains after the id in this paragraph (the lone
    # trailing-space run), keeping the paragraph mark itself intact.
    $newMatchEnd = $matchStart + $newId.Length
    $paraEnd = $p.Range.End - 1
    if ($paraEnd -gt $newMatchEnd) {
        $trailing = $d.Range($newMatchEnd, $paraEnd)
        $trailing.Delete()
    }
}

# Give the paragraph a border (5-twip gap on every side) and widen the left
# indent from 120 to 225 twips (6pt -> 11.25pt).
$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p.Range.ParagraphFormat.LeftIndent = 11.25

Write-Output "done"
